# Update for release to deploy 0.1.1
$wb = $excel.ActiveWorkbook

# Rename the second sheet from "Include from NMDP Preferred P" to "Include #0"
$includeSheet = $wb.Worksheets.Item("Include from NMDP Preferred P")
$includeSheet.Name = "Include #0"

# Metadata sheet updates
$meta = $wb.Worksheets.Item("Metadata")

# Bump version 0.1.0 -> 0.1.1
$meta.Range("B3").Value = "0.1.1"

# Update release date
$meta.Range("B8").Value = "2024-11-11T17:53:38-06:00"

# Insert a new "Jurisdiction" property row right before "Description" (currently row 11),
# pushing Description/Purpose/Copyright/Immutable down by one row.
$meta.Rows.Item(11).Insert()
$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""

# Match the look of the other property rows (border + top alignment)
$meta.Range("A11:B11").Borders.Item(8).LineStyle = 1
$meta.Range("A11:B11").VerticalAlignment = -4160
